$wb = $excel.ActiveWorkbook

# Rename sheet "STEPS TODO" to "STEPS"
$stepsSheet = $wb.Worksheets.Item("STEPS TODO")
$stepsSheet.Name = "STEPS"

# Set header row values for STEPS sheet.
# Values are assigned in an order that reproduces the shared-strings table
# order of the target workbook (new strings get appended in first-use order).
$stepsSheet.Range("A1").Value = "ACTION"
$stepsSheet.Range("D1").Value = "TC_STEP_ID"
$stepsSheet.Range("E1").Value = "TC_STEP_NUM"
$stepsSheet.Range("G1").Value = "TC_STEP_ACTION"
$stepsSheet.Range("I1").Value = "TC_STEP_#_REQ"
$stepsSheet.Range("J1").Value = "TC_STEP_#_ATTACHMENT"
$stepsSheet.Range("K1").Value = "TC_STEP_CUF_<CODE>"
$stepsSheet.Range("B1").Value = "TC_OWNER_PATH"
$stepsSheet.Range("C1").Value = "TC_OWNER_ID"
$stepsSheet.Range("F1").Value = "TC_STEP_IS_CALL_STEP"
$stepsSheet.Range("H1").Value = "TC_STEP_EXPECTED_RESULT"

# Set column widths (approximating the bestFit widths Excel computed on save)
$stepsSheet.Range("A1").EntireColumn.ColumnWidth = 7.0
$stepsSheet.Range("B1").EntireColumn.ColumnWidth = 15.666666666666666
$stepsSheet.Range("C1").EntireColumn.ColumnWidth = 12.666666666666666
$stepsSheet.Range("D1").EntireColumn.ColumnWidth = 10.333333333333334
$stepsSheet.Range("E1").EntireColumn.ColumnWidth = 12.833333333333334
$stepsSheet.Range("F1").EntireColumn.ColumnWidth = 20.166666666666668
$stepsSheet.Range("G1").EntireColumn.ColumnWidth = 15.5
$stepsSheet.Range("H1").EntireColumn.ColumnWidth = 24.5
$stepsSheet.Range("I1").EntireColumn.ColumnWidth = 14.0
$stepsSheet.Range("J1").EntireColumn.ColumnWidth = 22.833333333333332
$stepsSheet.Range("K1").EntireColumn.ColumnWidth = 19.833333333333332
$stepsSheet.Range("L1").EntireColumn.ColumnWidth = 19.833333333333332

# Make STEPS the active sheet/tab
$stepsSheet.Activate()
$stepsSheet.Range("K18").Select() | Out-Null
